$wb = $excel.ActiveWorkbook

# --- Add the new "system structure data" sheet after the existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "system structure data"

# --- Header row ---
$newSheet.Range("A1").Value = "System"
$newSheet.Range("B1").Value = "Structure"
$newSheet.Range("A1:B1").Font.Underline = 2

# --- Data rows. Written in an order that reproduces the author's original
#     shared-string insertion order (System, Structure, diesel engine system,
#     transmission system, [8,..], [2,3,4], fuel oil system, [(...)] ) while
#     still landing in the correct final cell positions. ---
$newSheet.Range("A3").Value = "diesel engine system"
$newSheet.Range("A2").Value = "transmission system"
$newSheet.Range("B3").Value = "[8, 9, 10, 11, 12, 6, 0, 1]"
$newSheet.Range("B2").Value = "[ 2, 3, 4]"

$newSheet.Range("A4").Value = "fuel oil system"
$newSheet.Range("B4").Value = "[([5,6], [5,6]), (7,7)]"

# --- Column widths ---
$newSheet.Range("A1:B1").ColumnWidth = 31.6

# --- Page setup ---
$newSheet.PageSetup.Orientation = 1

# --- Selection / active cell on the new sheet ---
$newSheet.Range("B5").Select() | Out-Null
